$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Area"
$ws.Range("B1").Value = "Considerations"
$ws.Range("C1").Value = "Normal Python"
$ws.Range("D1").Value = "SparkPython"
$ws.Range("E1").Value = "Related code file"

# Data rows
$ws.Range("A2").Value = "Data prep-processing"
$ws.Range("B2").Value = "Deduplication for the bag of words on the payment data level"

$ws.Range("A3").Value = "Data prep-processing"
$ws.Range("B3").Value = "Deduplication for the bag of words on the individual level"

$ws.Range("A4").Value = "Data prep-processing"
$ws.Range("B4").Value = "Deduplication of the failure bag on the payment data level"

$ws.Range("A5").Value = "Score Calculation"
$ws.Range("A6").Value = "ML model build"
$ws.Range("A7").Value = "ML model deployment"
$ws.Range("A8").Value = "Documentation "

# Header formatting: bold, light (theme background2) font color, accent1 (theme4) fill
$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.Font.ThemeColor = 4
$headerRange.Interior.ThemeColor = 5

# Column widths (best-fit approximations)
$ws.Columns.Item(1).ColumnWidth = 18.666666666666668
$ws.Columns.Item(2).ColumnWidth = 51
$ws.Columns.Item(3).ColumnWidth = 13
$ws.Columns.Item(4).ColumnWidth = 11
$ws.Columns.Item(5).ColumnWidth = 14.166666666666666

# Selection matching the target sheetView
[void]$ws.Range("A1:E8").Select()
